$wb = $excel.ActiveWorkbook

# Update "想去人数" (F2, F3) values on both the "展览" sheet and the
# "全部类型" sheet, which mirror the same event rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 214
    $ws.Range("F3").Value = 155
}
